$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Pepsi ---
$ws.Range("A2").Value = "Pepsi"
$ws.Range("C2").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/2/pepsico-24032024192227.png"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/2/pepsico-24032024192227.png")
$ws.Range("C2").Style = "Hyperlink"

# --- Row 3: Tra o Long Tea+ ---
$ws.Range("A3").Value = "Trà ô Long Tea+"
$ws.Range("C3").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/2/pepsico-24032024192227.png"
$ws.Hyperlinks.Add($ws.Range("C3"), "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/2/pepsico-24032024192227.png")
$ws.Range("C3").Style = "Hyperlink"

# --- Row 4: Mirinda ---
$ws.Range("A4").Value = "Mirinda"
$ws.Range("C4").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/5/pepsico-09042021143726.png"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/5/pepsico-09042021143726.png")
$ws.Range("C4").Style = "Hyperlink"

# --- Row 5: Sting ---
$ws.Range("A5").Value = "Sting"
$ws.Range("C5").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/12/pepsico-160920228953.png"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/12/pepsico-160920228953.png")
$ws.Range("C5").Style = "Hyperlink"

# --- Row 6: 7 up (keeps its existing hyperlink target, only the label/text changes) ---
$ws.Range("A6").Value = "7 up"
$ws.Range("C6").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/3/pepsico-24032024191926.png"

# --- Row 7: Revive ---
$ws.Range("A7").Value = "Revive"
$ws.Range("C7").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/13/pepsico-160920228937.png"

# --- Row 8: Rockstar ---
$ws.Range("A8").Value = "Rockstar"
$ws.Range("C8").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/14/pepsico-24032024192535.png"

# --- Row 9: Aquafina ---
$ws.Range("A9").Value = "Aquafina"
$ws.Range("C9").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/8/pepsico-29082022105643.png"

# --- Row 10: Twister ---
$ws.Range("A10").Value = "Twister"
$ws.Range("C10").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/9/pepsico-09042021143842.png"

# --- Row 11: Lipton (new row, no hyperlink/style) ---
$ws.Range("A11").Value = "Lipton"
$ws.Range("C11").Value = "https://cdn.tgdd.vn/bachhoaxanh/shopinshop/8/10/pepsico-31082022172335.png"

$ws.Range("C11").Select()
